# Hybrid_Dewatering_BOM.xlsx -- "Updating of the total Cost of the BOM"
#
# Adds a TOTAL row (row 19) to the PLC_BOM sheet with SUM formulas for the
# Low/High unit-price and Low/High total columns, plus a standalone
# "TOTAL Cost Approx" callout at I12:J12. Also refreshes the view state
# (zoom/selection/active sheet) and auto-fits columns on every sheet, the
# way Excel would after a user reviewed the whole workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PLC_BOM")

# --- New TOTAL row (row 19) --------------------------------------------
$ws.Range("A19").Value = "TOTAL"
$ws.Range("A19").HorizontalAlignment = -4152   # xlRight

$ws.Range("C19").Formula = "=SUM(C2:C18)"
$ws.Range("F19").Formula = "=SUM(F2:F18)"
$ws.Range("G19").Formula = "=SUM(G2:G18)"

# --- New standalone callout cell (I12 / J12) ---------------------------
$ws.Range("I12").Value = "TOTAL Cost Approx"
$ws.Range("J12").Value = 2000000

# --- Recalculate so cached formula results are up to date --------------
$excel.Calculate()

# --- Autofit columns on every sheet (whole-workbook review pass) -------
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Select()
    $sheet.Cells.EntireColumn.AutoFit()
}

# --- Restore per-sheet selections / zoom as left by the author ---------
$plc = $wb.Worksheets.Item("PLC_BOM")
$plc.Range("I18").Select()
$excel.ActiveWindow.Zoom = 118

$solar = $wb.Worksheets.Item("Solar_Panel_Calc")
$solar.Cells.Select()
$excel.ActiveWindow.Zoom = 82

$installed = $wb.Worksheets.Item("Installed_Plant")
$installed.Range("C2").Select()

$vfd = $wb.Worksheets.Item("VFD_and_Pump")
$vfd.Cells.Select()

$notes = $wb.Worksheets.Item("Notes")
$notes.Range("B17").Select()

$totals = $wb.Worksheets.Item("Totals_Summary")
$totals.Range("B10").Select()
$totals.Activate()

Write-Host "PLC_BOM TOTAL row + TOTAL Cost Approx callout added."
